$d = $word.ActiveDocument
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "SEQ Table") {
        $f.Delete()
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Table") {
        $insertPos = $p.Range.Start + 6
        $r = $d.Range($insertPos, $insertPos)
        $newField = $d.Fields.Add($r, 12, " SEQ Table \* ARABIC ", $false)
    }
}
foreach ($f2 in $d.Fields) {
    if ($f2.Code.Text -match "SEQ Table") {
        $f2.Result.Bold = 1
        Write-Host "set bold via fresh field object result"
    }
}
